# Apply cryptos.xlsx price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.629.63"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "3.785.48"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.58"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.66"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "3.786.03"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.45"
$ws.Range("D15").Value = "4.419.99"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "3.796.71"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.49"
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("D18").Value = "67.561.34"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.19"
$ws.Range("E21").Value = "  -4.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "456.45"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000155"
$ws.Range("E24").Value = "  +7.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.49"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.91"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.85"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.23"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").Value = "3.737.31"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.26"
$ws.Range("E44").Value = "  +5.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.300"
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.07"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.36"
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "149.03"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("E49").Value = "  -4.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "389.59"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.756.22"
$ws.Range("E51").Value = "  +2.23%  "
